$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    # Force text storage so numeric/percentage-looking strings are kept as
    # literal text, matching the original inlineStr cell typing instead of
    # being auto-converted to a number/percentage by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - BNB
Set-TextValue "D2" "294.63"
Set-TextValue "E2" "-4.68%"

# Row 3 - OKB
Set-TextValue "D3" "40.14"
Set-TextValue "E3" "-2.99%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.021"
Set-TextValue "E4" "-3.73%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07381"

# Row 6 - GateToken
Set-TextValue "D6" "4.306"
Set-TextValue "E6" "-0.26%"

# Row 7 - FTXToken
Set-TextValue "D7" "1.562"
Set-TextValue "E7" "-4.79%"

# Row 8 - MXToken
Set-TextValue "D8" "0.9225"

# Row 9 - LiechtensteinCryptoassetsExchange
Set-TextValue "D9" "0.1189"
Set-TextValue "E9" "-3.83%"

# Row 10 - WazirX
Set-TextValue "D10" "0.1775"
Set-TextValue "E10" "-2.63%"

# Row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.08743"
Set-TextValue "E11" "-4.62%"

# Row 12 - BitrueCoin
Set-TextValue "D12" "0.04168"
Set-TextValue "E12" "-0.27%"

# Row 13 - BitMartToken
Set-TextValue "D13" "0.1053"
Set-TextValue "E13" "0.15%"

# Row 14 - BitForexToken
Set-TextValue "D14" "0.001277"
Set-TextValue "E14" "1.95%"

# Row 15 - was CoinExToken, now TigerCash
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005779"
Set-TextValue "E15" "-1.19%"

# Row 16 - was TigerCash, now LEO
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.382"
Set-TextValue "E16" "0.98%"

# Row 17 - was LEO, now BTSEToken
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D17" "2.399"
Set-TextValue "E17" "-1.07%"

# Row 18 - was BTSEToken, now BitpandaEcosystemToken
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D18" "0.3296"
Set-TextValue "E18" "-0.73%"

# Row 19 - was BitpandaEcosystemToken, now MCDex
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D19" "7.550"
Set-TextValue "E19" "2.31%"

# Row 20 - was MCDex, now ProBitToken
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1344"
Set-TextValue "E20" "-4.17%"

# Row 21 - was ProBitToken, now ZBToken
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D21" "0.2808"
Set-TextValue "E21" "-0.51%"

# Row 22 - was ZBToken, now CoinExToken
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D22" "0.03808"
Set-TextValue "E22" "-5.15%"

# Row 23 - BitKan
Set-TextValue "D23" "0.001285"
Set-TextValue "E23" "1.44%"

# Row 24 - HotbitToken
Set-TextValue "D24" "0.003904"
Set-TextValue "E24" "-4.49%"

# Row 25 - NitroEx
Set-TextValue "D25" "0.0001293"
Set-TextValue "E25" "-0.64%"

# Row 26 - UpBots
Set-TextValue "D26" "0.0003731"
Set-TextValue "E26" "-95.03%"

# Row 38 - One
Set-TextValue "D38" "0.02313"
Set-TextValue "E38" "-8.97%"

# Row 39 - IDEX
Set-TextValue "D39" "0.05036"
Set-TextValue "E39" "-5.81%"

# Row 40 - KickToken
Set-TextValue "D40" "0.007738"
Set-TextValue "E40" "-1.38%"

# Row 41 - CEJI
Set-TextValue "E41" "137.24%"

# Row 42 - BKEXToken
Set-TextValue "D42" "0.1276"
Set-TextValue "E42" "-2.92%"

# Row 43 - Dexo
Set-TextValue "D43" "0.007408"
Set-TextValue "E43" "10.89%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.006981"
Set-TextValue "E44" "-5.40%"

# Row 45 - PooCoin
Set-TextValue "D45" "0.3196"
Set-TextValue "E45" "4.16%"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00006465"
Set-TextValue "E46" "-4.55%"

# Row 47 - Kangarootoken
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.12%"

# Row 48 - BOLO
Set-TextValue "E48" "11.82%"

# Row 49 - CoinbaseStockToken
Set-TextValue "D49" "0.004209"
Set-TextValue "E49" "35.63%"

# Row 50 - CryptobidCoin
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.12%"

# Row 51 - SpecialPowerGold
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.12%"
